# Insert one new weekly price-report row above row 507 for
# "Hortaliza, Vega Central Mapocho de Santiago - Albahaca".
# Excel shifts the existing rows 507:618 down to 508:619 and the
# newly-blank row 507 is then filled with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(507).Insert()

$ws.Range("A507").Value = 9
$ws.Range("B507").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C507").Value = "Metropolitana"
$ws.Range("D507").Value = 45173
$ws.Range("E507").Value = 13
$ws.Range("F507").Value = 100112052
$ws.Range("G507").Value = "Albahaca"
$ws.Range("H507").Value = "Sin especificar"
$ws.Range("I507").Value = "Primera"
$ws.Range("J507").Value = 430
$ws.Range("K507").Value = 3500
$ws.Range("L507").Value = 4000
$ws.Range("M507").Value = 3750
$ws.Range("N507").Value = "$/paquete"
$ws.Range("O507").Value = "Región de Arica y Parinacota"
$ws.Range("P507").Value = 3750
$ws.Range("Q507").Value = 1
$ws.Range("R507").Value = "Hortaliza"
